$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 12.785
$ws.Range("E6").Value = 12.711
$ws.Range("E7").Value = 13.201
$ws.Range("E8").Value = 12.668
$ws.Range("E16").Value = 12.668
$ws.Range("E20").Value = 12.725
$ws.Range("E21").Value = 13.33
